# Actualización automática: 2025-03-21 10:50:08
# Appends two new inventory rows (CRM / ansible_test hosts) to the sheet,
# mirroring the pattern of the existing rows, and updates the active
# selection to match the new last-used cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 130: 10.181.5.214 / contactcrmdbprd
$ws.Range("A130").Value = "10.181.5.214"
$ws.Range("B130").Value = "CRM"
$ws.Range("C130").Value = "10.181.5.214"
$ws.Range("D130").Value = "ansible_test"
$ws.Range("E130").Value = "contactcrmdbprd"

# Row 131: 10.181.5.219 / crmdb4
$ws.Range("A131").Value = "10.181.5.219"
$ws.Range("B131").Value = "CRM"
$ws.Range("C131").Value = "10.181.5.219"
$ws.Range("D131").Value = "ansible_test"
$ws.Range("E131").Value = "crmdb4"

# Match the author's final selection/cursor position.
[void]$ws.Range("E131").Select()
